$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the bold/bordered "A" column style from the last existing data row
# (row 16) down onto the three newly appended rows (17-19) before filling
# in their values, so the new rows render consistently with the rest of
# the table.
$ws.Cells.Item(16, 1).Copy($ws.Cells.Item(17, 1))
$ws.Cells.Item(16, 1).Copy($ws.Cells.Item(18, 1))
$ws.Cells.Item(16, 1).Copy($ws.Cells.Item(19, 1))

# Refresh the causality-test statistics for the BXD strain panel: update
# existing rows 2-16 and populate the newly added rows 17-19 (BXD6, BXD16,
# BXD9) with their likelihood/probability results.

# Row 2: BXD1
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "BXD1"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 3.1855
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 38.875
$ws.Cells.Item(2, 7).Value = 1245.109375
$ws.Cells.Item(2, 8).Value = 2.99170483610721
$ws.Cells.Item(2, 9).Value = 0.2070144549413103
$ws.Cells.Item(2, 10).Value = 0.006162279433048035
$ws.Cells.Item(2, 11).Value = 0.8007830544265581
$ws.Cells.Item(2, 12).Value = 0.002467324473313082

# Row 3: BXD5
$ws.Cells.Item(3, 1).Value = 27
$ws.Cells.Item(3, 2).Value = "BXD5"
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 3.0645
$ws.Cells.Item(3, 5).Value = 83
$ws.Cells.Item(3, 6).Value = 38.875
$ws.Cells.Item(3, 7).Value = 1245.109375
$ws.Cells.Item(3, 8).Value = 3.277132813667195
$ws.Cells.Item(3, 9).Value = 0.2070144549413103
$ws.Cells.Item(3, 10).Value = 0.005173048256839504
$ws.Cells.Item(3, 11).Value = 0.7861111914692198
$ws.Cells.Item(3, 12).Value = 0.002033295564355937

# Row 4: BXD31
$ws.Cells.Item(4, 1).Value = 14
$ws.Cells.Item(4, 2).Value = "BXD31"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 2.931
$ws.Cells.Item(4, 5).Value = 20
$ws.Cells.Item(4, 6).Value = 38.875
$ws.Cells.Item(4, 7).Value = 1245.109375
$ws.Cells.Item(4, 8).Value = 3.060482662025279
$ws.Cells.Item(4, 9).Value = 0.2070144549413103
$ws.Cells.Item(4, 10).Value = 0.00979881462615984
$ws.Cells.Item(4, 11).Value = 0.8420217944128214
$ws.Cells.Item(4, 12).Value = 0.004125407737318855

# Row 5: BXD23
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "BXD23"
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 2.6615
$ws.Cells.Item(5, 5).Value = 20
$ws.Cells.Item(5, 6).Value = 38.875
$ws.Cells.Item(5, 7).Value = 1245.109375
$ws.Cells.Item(5, 8).Value = 3.060482662025279
$ws.Cells.Item(5, 9).Value = 0.2070144549413103
$ws.Cells.Item(5, 10).Value = 0.00979881462615984
$ws.Cells.Item(5, 11).Value = 0.5969401540119526
$ws.Cells.Item(5, 12).Value = 0.002924652956037214

# Row 6: BXD2
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = "BXD2"
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 3.2065
$ws.Cells.Item(6, 5).Value = 10
$ws.Cells.Item(6, 6).Value = 38.875
$ws.Cells.Item(6, 7).Value = 1245.109375
$ws.Cells.Item(6, 8).Value = 3.026093749066245
$ws.Cells.Item(6, 9).Value = 0.2070144549413103
$ws.Cells.Item(6, 10).Value = 0.008089049222233136
$ws.Cells.Item(6, 11).Value = 0.8105321716427392
$ws.Cells.Item(6, 12).Value = 0.003278217316310817

# Row 7: BXD21
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = "BXD21"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 2.374
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 38.875
$ws.Cells.Item(7, 7).Value = 1245.109375
$ws.Cells.Item(7, 8).Value = 3.026093749066245
$ws.Cells.Item(7, 9).Value = 0.2070144549413103
$ws.Cells.Item(7, 10).Value = 0.008089049222233136
$ws.Cells.Item(7, 11).Value = 0.3139567474823901
$ws.Cells.Item(7, 12).Value = 0.001269805792018636

# Row 8: BXD14
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 2).Value = "BXD14"
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 2.6955
$ws.Cells.Item(8, 5).Value = 89
$ws.Cells.Item(8, 6).Value = 38.875
$ws.Cells.Item(8, 7).Value = 1245.109375
$ws.Cells.Item(8, 8).Value = 3.297766161442615
$ws.Cells.Item(8, 9).Value = 0.2070144549413103
$ws.Cells.Item(8, 10).Value = 0.004122144904925831
$ws.Cells.Item(8, 11).Value = 0.3651159762517838
$ws.Cells.Item(8, 12).Value = 0.0007525304806066557

# Row 9: BXD13
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 2).Value = "BXD13"
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 3.9685
$ws.Cells.Item(9, 5).Value = 79
$ws.Cells.Item(9, 6).Value = 38.875
$ws.Cells.Item(9, 7).Value = 1245.109375
$ws.Cells.Item(9, 8).Value = 3.263377248483581
$ws.Cells.Item(9, 9).Value = 0.2070144549413103
$ws.Cells.Item(9, 10).Value = 0.005922695085047902
$ws.Cells.Item(9, 11).Value = 0.2638609740989317
$ws.Cells.Item(9, 12).Value = 0.0007813840472158472

# Row 10: BXD19
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "BXD19"
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 2.905
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 45.3
$ws.Cells.Item(10, 7).Value = 1960.41
$ws.Cells.Item(10, 8).Value = 2.99170483610721
$ws.Cells.Item(10, 9).Value = 0.2070144549413103
$ws.Cells.Item(10, 10).Value = 0.005338685304373844
$ws.Cells.Item(10, 11).Value = 0.8610413921187994
$ws.Cells.Item(10, 12).Value = 0.002298414513281115

# Row 11: BXD15
$ws.Cells.Item(11, 1).Value = 5
$ws.Cells.Item(11, 2).Value = "BXD15"
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 3.3835
$ws.Cells.Item(11, 5).Value = 100
$ws.Cells.Item(11, 6).Value = 45.3
$ws.Cells.Item(11, 7).Value = 1960.41
$ws.Cells.Item(11, 8).Value = 3.335593965697553
$ws.Cells.Item(11, 9).Value = 0.2070144549413103
$ws.Cells.Item(11, 10).Value = 0.004200626402036855
$ws.Cells.Item(11, 11).Value = 0.871971737309845
$ws.Cells.Item(11, 12).Value = 0.00183141375078684

# Row 12: BXD8
$ws.Cells.Item(12, 1).Value = 35
$ws.Cells.Item(12, 2).Value = "BXD8"
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 3.353
$ws.Cells.Item(12, 5).Value = 100
$ws.Cells.Item(12, 6).Value = 45.3
$ws.Cells.Item(12, 7).Value = 1960.41
$ws.Cells.Item(12, 8).Value = 3.335593965697553
$ws.Cells.Item(12, 9).Value = 0.2070144549413103
$ws.Cells.Item(12, 10).Value = 0.004200626402036855
$ws.Cells.Item(12, 11).Value = 0.8761771652902894
$ws.Cells.Item(12, 12).Value = 0.0018402464666901

# Row 13: BXD24
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "BXD24"
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 2.995
$ws.Cells.Item(13, 5).Value = 100
$ws.Cells.Item(13, 6).Value = 45.3
$ws.Cells.Item(13, 7).Value = 1960.41
$ws.Cells.Item(13, 8).Value = 3.335593965697553
$ws.Cells.Item(13, 9).Value = 0.2070144549413103
$ws.Cells.Item(13, 10).Value = 0.004200626402036855
$ws.Cells.Item(13, 11).Value = 0.6625633254156288
$ws.Cells.Item(13, 12).Value = 0.001391590498881113

# Row 14: BXD28
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "BXD28"
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 4.2805
$ws.Cells.Item(14, 5).Value = 66
$ws.Cells.Item(14, 6).Value = 45.3
$ws.Cells.Item(14, 7).Value = 1960.41
$ws.Cells.Item(14, 8).Value = 3.218671661636836
$ws.Cells.Item(14, 9).Value = 0.2070144549413103
$ws.Cells.Item(14, 10).Value = 0.008077451998695988
$ws.Cells.Item(14, 11).Value = 0.05757624069766825
$ws.Cells.Item(14, 12).Value = 0.0002325346602503909

# Row 15: BXD32
$ws.Cells.Item(15, 1).Value = 15
$ws.Cells.Item(15, 2).Value = "BXD32"
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 3.072
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 45.3
$ws.Cells.Item(15, 7).Value = 1960.41
$ws.Cells.Item(15, 8).Value = 2.99170483610721
$ws.Cells.Item(15, 9).Value = 0.2070144549413103
$ws.Cells.Item(15, 10).Value = 0.005338685304373844
$ws.Cells.Item(15, 11).Value = 0.8632703772916952
$ws.Cells.Item(15, 12).Value = 0.002304364438474218

# Row 16: BXD11
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "BXD11"
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 3.607
$ws.Cells.Item(16, 5).Value = 12
$ws.Cells.Item(16, 6).Value = 45.3
$ws.Cells.Item(16, 7).Value = 1960.41
$ws.Cells.Item(16, 8).Value = 3.032971531658052
$ws.Cells.Item(16, 9).Value = 0.2070144549413103
$ws.Cells.Item(16, 10).Value = 0.006790613889094571
$ws.Cells.Item(16, 11).Value = 0.3956147791747662
$ws.Cells.Item(16, 12).Value = 0.001343233607097624

# Row 17: BXD6
$ws.Cells.Item(17, 1).Value = 29
$ws.Cells.Item(17, 2).Value = "BXD6"
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 3.387
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 45.3
$ws.Cells.Item(17, 7).Value = 1960.41
$ws.Cells.Item(17, 8).Value = 2.99170483610721
$ws.Cells.Item(17, 9).Value = 0.2070144549413103
$ws.Cells.Item(17, 10).Value = 0.005338685304373844
$ws.Cells.Item(17, 11).Value = 0.6011779633041984
$ws.Cells.Item(17, 12).Value = 0.001604749979002761

# Row 18: BXD16
$ws.Cells.Item(18, 1).Value = 6
$ws.Cells.Item(18, 2).Value = "BXD16"
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 2.4065
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 45.3
$ws.Cells.Item(18, 7).Value = 1960.41
$ws.Cells.Item(18, 8).Value = 2.99170483610721
$ws.Cells.Item(18, 9).Value = 0.2070144549413103
$ws.Cells.Item(18, 10).Value = 0.005338685304373844
$ws.Cells.Item(18, 11).Value = 0.3834266793074632
$ws.Cells.Item(18, 12).Value = 0.001023497189061808

# Row 19: BXD9
$ws.Cells.Item(19, 1).Value = 38
$ws.Cells.Item(19, 2).Value = "BXD9"
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 3.0015
$ws.Cells.Item(19, 5).Value = 75
$ws.Cells.Item(19, 6).Value = 45.3
$ws.Cells.Item(19, 7).Value = 1960.41
$ws.Cells.Item(19, 8).Value = 3.249621683299968
$ws.Cells.Item(19, 9).Value = 0.2070144549413103
$ws.Cells.Item(19, 10).Value = 0.007195000363549942
$ws.Cells.Item(19, 11).Value = 0.7556696055810315
$ws.Cells.Item(19, 12).Value = 0.002718521543439581
